$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 143, shifting existing rows 143:229 down to 144:230.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new record.
$ws.Cells.Item(143, 1).Value = 3
$ws.Cells.Item(143, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(143, 3).Value = "Coquimbo"
$ws.Cells.Item(143, 4).Value = 44830
$ws.Cells.Item(143, 5).Value = 5
$ws.Cells.Item(143, 6).Value = "Fruta"
$ws.Cells.Item(143, 7).Value = 100101
$ws.Cells.Item(143, 8).Value = "Berries"
$ws.Cells.Item(143, 9).Value = 100101001
$ws.Cells.Item(143, 10).Value = "Arándano (blue)"
$ws.Cells.Item(143, 11).Value = "Sin especificar"
$ws.Cells.Item(143, 12).Value = "Primera"
$ws.Cells.Item(143, 13).Value = 56
$ws.Cells.Item(143, 14).Value = 11000
$ws.Cells.Item(143, 15).Value = 11000
$ws.Cells.Item(143, 16).Value = 11000
$ws.Cells.Item(143, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(143, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(143, 19).Value = 7333
$ws.Cells.Item(143, 20).Value = 1.5
